# Updates 2-gram frequency tables on both sheets (negative/positive) to
# remove n-grams that are shared between the two sheets, per commit
# "Results removed shared words". Rows 2-65 (A: n-gram text, B: count, C: proportion)
# are replaced with the recomputed top-64 lists for each sheet.

$wb = $excel.ActiveWorkbook

$negativeData = @(
  @('(''waste'', ''time'')', 103, 0.04561558901682905),
  @('(''really'', ''bad'')', 70, 0.03100088573959256),
  @('(''pretty'', ''much'')', 61, 0.02701505757307352),
  @('(''bad'', ''acting'')', 60, 0.02657218777679362),
  @('(''bad'', ''movie'')', 60, 0.02657218777679362),
  @('(''movie'', ''bad'')', 54, 0.02391496899911426),
  @('(''one'', ''worst'')', 53, 0.02347209920283437),
  @('(''movie'', ''even'')', 43, 0.01904340124003543),
  @('(''bad'', ''guys'')', 42, 0.01860053144375554),
  @('(''first'', ''minutes'')', 42, 0.01860053144375554),
  @('(''movies'', ''ever'')', 41, 0.01815766164747564),
  @('(''worst'', ''movie'')', 41, 0.01815766164747564),
  @('(''even'', ''worse'')', 41, 0.01815766164747564),
  @('(''looked'', ''like'')', 40, 0.01771479185119575),
  @('(''end'', ''movie'')', 39, 0.01727192205491586),
  @('(''character'', ''development'')', 38, 0.01682905225863596),
  @('(''make'', ''sense'')', 37, 0.01638618246235607),
  @('(''b'', ''movie'')', 36, 0.01594331266607617),
  @('(''writer'', ''director'')', 36, 0.01594331266607617),
  @('(''something'', ''like'')', 35, 0.01550044286979628),
  @('(''acting'', ''bad'')', 34, 0.01505757307351639),
  @('(''thing'', ''movie'')', 34, 0.01505757307351639),
  @('(''main'', ''characters'')', 34, 0.01505757307351639),
  @('(''film'', ''ever'')', 34, 0.01505757307351639),
  @('(''film'', ''making'')', 34, 0.01505757307351639),
  @('(''good'', ''thing'')', 34, 0.01505757307351639),
  @('(''make'', ''film'')', 34, 0.01505757307351639),
  @('(''felt'', ''like'')', 33, 0.01461470327723649),
  @('(''horror'', ''movies'')', 33, 0.01461470327723649),
  @('(''worst'', ''movies'')', 33, 0.01461470327723649),
  @('(''movie'', ''could'')', 33, 0.01461470327723649),
  @('(''movie'', ''movie'')', 32, 0.0141718334809566),
  @('(''two'', ''hours'')', 31, 0.0137289636846767),
  @('(''end'', ''film'')', 31, 0.0137289636846767),
  @('(''time'', ''money'')', 31, 0.0137289636846767),
  @('(''one'', ''point'')', 31, 0.0137289636846767),
  @('(''bad'', ''film'')', 31, 0.0137289636846767),
  @('(''first'', ''movie'')', 30, 0.01328609388839681),
  @('(''whole'', ''thing'')', 30, 0.01328609388839681),
  @('(''anything'', ''else'')', 29, 0.01284322409211692),
  @('(''first'', ''place'')', 29, 0.01284322409211692),
  @('(''yet'', ''another'')', 28, 0.01240035429583702),
  @('(''piece'', ''crap'')', 28, 0.01240035429583702),
  @('(''made'', ''tv'')', 28, 0.01240035429583702),
  @('(''camera'', ''work'')', 28, 0.01240035429583702),
  @('(''really'', ''really'')', 28, 0.01240035429583702),
  @('(''video'', ''store'')', 28, 0.01240035429583702),
  @('(''film'', ''could'')', 27, 0.01195748449955713),
  @('(''feels'', ''like'')', 27, 0.01195748449955713),
  @('(''whole'', ''film'')', 27, 0.01195748449955713),
  @('(''feel'', ''like'')', 27, 0.01195748449955713),
  @('(''seems'', ''like'')', 27, 0.01195748449955713),
  @('(''like'', ''watching'')', 27, 0.01195748449955713),
  @('(''would'', ''make'')', 26, 0.01151461470327724),
  @('(''one'', ''two'')', 26, 0.01151461470327724),
  @('(''read'', ''book'')', 26, 0.01151461470327724),
  @('(''production'', ''values'')', 26, 0.01151461470327724),
  @('(''seemed'', ''like'')', 26, 0.01151461470327724),
  @('(''watching'', ''film'')', 26, 0.01151461470327724),
  @('(''part'', ''movie'')', 25, 0.01107174490699734),
  @('(''bad'', ''movies'')', 25, 0.01107174490699734),
  @('(''first'', ''one'')', 25, 0.01107174490699734),
  @('(''absolutely'', ''nothing'')', 25, 0.01107174490699734),
  @('(''nothing'', ''new'')', 25, 0.01107174490699734)
)

$positiveData = @(
  @('(''one'', ''best'')', 130, 0.04684684684684685),
  @('(''new'', ''york'')', 70, 0.03434739941118744),
  @('(''must'', ''see'')', 69, 0.03385672227674191),
  @('(''great'', ''movie'')', 61, 0.02993130520117762),
  @('(''highly'', ''recommend'')', 51, 0.02502453385672228),
  @('(''love'', ''story'')', 47, 0.02306182531894014),
  @('(''years'', ''later'')', 45, 0.02208047105004907),
  @('(''one'', ''favorite'')', 41, 0.02011776251226693),
  @('(''black'', ''white'')', 41, 0.02011776251226693),
  @('(''great'', ''film'')', 41, 0.02011776251226693),
  @('(''love'', ''movie'')', 38, 0.01864573110893033),
  @('(''well'', ''done'')', 37, 0.01815505397448479),
  @('(''tv'', ''series'')', 37, 0.01815505397448479),
  @('(''great'', ''job'')', 36, 0.01766437684003926),
  @('(''tony'', ''hawk'')', 35, 0.01717369970559372),
  @('(''world'', ''war'')', 35, 0.01717369970559372),
  @('(''well'', ''worth'')', 35, 0.01717369970559372),
  @('(''good'', ''job'')', 34, 0.01668302257114818),
  @('(''true'', ''story'')', 34, 0.01668302257114818),
  @('(''recommend'', ''movie'')', 33, 0.01619234543670265),
  @('(''film'', ''festival'')', 33, 0.01619234543670265),
  @('(''first'', ''saw'')', 31, 0.01521099116781158),
  @('(''film'', ''great'')', 31, 0.01521099116781158),
  @('(''young'', ''man'')', 30, 0.01472031403336605),
  @('(''every'', ''time'')', 30, 0.01472031403336605),
  @('(''films'', ''like'')', 29, 0.01422963689892051),
  @('(''film'', ''really'')', 29, 0.01422963689892051),
  @('(''movie'', ''great'')', 28, 0.01373895976447497),
  @('(''movie'', ''seen'')', 28, 0.01373895976447497),
  @('(''supporting'', ''cast'')', 27, 0.01324828263002944),
  @('(''watch'', ''film'')', 27, 0.01324828263002944),
  @('(''th'', ''century'')', 27, 0.01324828263002944),
  @('(''fun'', ''watch'')', 26, 0.01275760549558391),
  @('(''movie'', ''also'')', 25, 0.01226692836113837),
  @('(''martial'', ''arts'')', 25, 0.01226692836113837),
  @('(''one'', ''greatest'')', 25, 0.01226692836113837),
  @('(''kung'', ''fu'')', 25, 0.01226692836113837),
  @('(''highly'', ''recommended'')', 25, 0.01226692836113837),
  @('(''top'', ''notch'')', 24, 0.01177625122669284),
  @('(''big'', ''screen'')', 24, 0.01177625122669284),
  @('(''give'', ''movie'')', 24, 0.01177625122669284),
  @('(''quite'', ''good'')', 24, 0.01177625122669284),
  @('(''also'', ''good'')', 24, 0.01177625122669284),
  @('(''watched'', ''movie'')', 23, 0.0112855740922473),
  @('(''really'', ''like'')', 23, 0.0112855740922473),
  @('(''great'', ''acting'')', 23, 0.0112855740922473),
  @('(''put'', ''together'')', 23, 0.0112855740922473),
  @('(''along'', ''way'')', 23, 0.0112855740922473),
  @('(''hong'', ''kong'')', 23, 0.0112855740922473),
  @('(''film'', ''well'')', 23, 0.0112855740922473),
  @('(''seen'', ''movie'')', 22, 0.01079489695780177),
  @('(''would'', ''recommend'')', 22, 0.01079489695780177),
  @('(''film'', ''noir'')', 22, 0.01079489695780177),
  @('(''based'', ''true'')', 22, 0.01079489695780177),
  @('(''pro'', ''skater'')', 22, 0.01079489695780177),
  @('(''hawk'', ''pro'')', 22, 0.01079489695780177),
  @('(''saw'', ''film'')', 22, 0.01079489695780177),
  @('(''throughout'', ''film'')', 22, 0.01079489695780177),
  @('(''would'', ''never'')', 22, 0.01079489695780177),
  @('(''feel'', ''good'')', 22, 0.01079489695780177),
  @('(''sense'', ''humor'')', 22, 0.01079489695780177),
  @('(''old'', ''man'')', 22, 0.01079489695780177),
  @('(''well'', ''acted'')', 21, 0.01030421982335623),
  @('(''horror'', ''films'')', 21, 0.01030421982335623)
)

$wsNegative = $wb.Worksheets.Item("negative")
for ($i = 0; $i -lt $negativeData.Count; $i++) {
    $row = $negativeData[$i]
    $r = $i + 2
    $wsNegative.Cells.Item($r, 1).Value = $row[0]
    $wsNegative.Cells.Item($r, 2).Value = $row[1]
    $wsNegative.Cells.Item($r, 3).Value = $row[2]
}

$wsPositive = $wb.Worksheets.Item("positive")
for ($i = 0; $i -lt $positiveData.Count; $i++) {
    $row = $positiveData[$i]
    $r = $i + 2
    $wsPositive.Cells.Item($r, 1).Value = $row[0]
    $wsPositive.Cells.Item($r, 2).Value = $row[1]
    $wsPositive.Cells.Item($r, 3).Value = $row[2]
}
